# Allow alt_cells for the relative formula type.
# Adds two new rows of relative-sample data (rows 4 and 5) and updates
# the active selection on the "Relative Samples" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relative Samples")

# Rows 4 and 5 were previously empty (data resumed at row 6), so this
# just fills in the gap with two more sample rows, shaped like the
# existing ones in columns A:C.
$ws.Range("A4").Value = 3005
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 2000

$ws.Range("A5").Value = 8888
$ws.Range("B5").Value = 1000
$ws.Range("C5").Value = 2000

# Update the active selection on the sheet.
$ws.Activate()
$ws.Range("F13").Select()
